$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.673.75"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.696.75"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.77"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3716"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.93"
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3399"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.202"
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07429"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.292"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.81"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.949"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "1.695.89"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001114"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06695"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "82.88"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.09"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.297"
$ws.Range("E22").Value = "  +3.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.88"
$ws.Range("E23").Value = "  +6.31%  "
$ws.Range("D24").Value = "24.705.85"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.456"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.745"
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.12"
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.50"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "131.15"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.883.94"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.220"
$ws.Range("E31").Value = "  +25.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.666"
$ws.Range("E32").Value = "  +6.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.210"
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08672"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.752"
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.46"
$ws.Range("E36").Value = "  +9.09%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.503"
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06583"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.976"
$ws.Range("E39").Value = "  +4.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02387"
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2202"
$ws.Range("E41").Value = "  +5.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.253"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6357"
$ws.Range("E43").Value = "  +4.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.65"
$ws.Range("E45").Value = "  +5.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6052"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.809"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.098"
$ws.Range("E48").Value = "  +4.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.32"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07230"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.95"
$ws.Range("E51").Value = "  +4.33%  "
